$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text so numeric-looking strings
# (e.g. trailing zeros, percent signs) are preserved verbatim.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.316.06'
$ws.Range('E2').Value = '  +2.03%  '
$ws.Range('D3').Value = '2.092.21'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.48%  '
$ws.Range('D5').Value = '343.33'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('D7').Value = '0.5210'
$ws.Range('E7').Value = '  +1.07%  '
$ws.Range('D8').Value = '0.4405'
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').Value = '54.30'
$ws.Range('E9').Value = '  +2.34%  '
$ws.Range('D10').Value = '0.09330'
$ws.Range('E10').Value = '  +1.22%  '
$ws.Range('D11').Value = '1.170'
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').Value = '24.79'
$ws.Range('E12').Value = '  -0.66%  '
$ws.Range('D13').Value = '8.689'
$ws.Range('E13').Value = '  +5.92%  '
$ws.Range('D14').Value = '2.115.16'
$ws.Range('E14').Value = '  +0.95%  '
$ws.Range('D15').Value = '6.908'
$ws.Range('E15').Value = '  +2.18%  '
$ws.Range('D16').Value = '101.21'
$ws.Range('E16').Value = '  +1.69%  '
$ws.Range('D17').Value = '0.00001157'
$ws.Range('E17').Value = '  +0.51%  '
$ws.Range('D18').Value = '1.004'
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('D19').Value = '21.20'
$ws.Range('E19').Value = '  +2.35%  '
$ws.Range('D20').Value = '0.06684'
$ws.Range('E20').Value = '  +0.77%  '
$ws.Range('D21').Value = '6.373'
$ws.Range('E21').Value = '  +2.94%  '
$ws.Range('D22').Value = '1.004'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').Value = '30.298.55'
$ws.Range('E23').Value = '  +1.85%  '
$ws.Range('D24').Value = '12.53'
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('D25').Value = '2.294'
$ws.Range('E25').Value = '  -0.90%  '
$ws.Range('D26').Value = '21.72'
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('D27').Value = '162.05'
$ws.Range('E27').Value = '  -0.35%  '
$ws.Range('D28').Value = '2.521'
$ws.Range('E28').Value = '  -0.42%  '
$ws.Range('D29').Value = '132.93'
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('D30').Value = '1.132'
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').Value = '1.671'
$ws.Range('E31').Value = '  +1.28%  '
$ws.Range('D32').Value = '0.1048'
$ws.Range('E32').Value = '  -0.12%  '
$ws.Range('D33').Value = '6.221'
$ws.Range('E33').Value = '  +0.91%  '
$ws.Range('D34').Value = '6.697'
$ws.Range('E34').Value = '  +10.58%  '
$ws.Range('D35').Value = '3.847'
$ws.Range('E35').Value = '  -2.78%  '
$ws.Range('D36').Value = '10.21'
$ws.Range('E36').Value = '  -1.57%  '
$ws.Range('D37').Value = '0.02629'
$ws.Range('E37').Value = '  +2.56%  '
$ws.Range('D38').Value = '0.06770'
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = '0.6989'
$ws.Range('E39').Value = '  +1.78%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '1.348'
$ws.Range('E40').Value = '  +4.22%  '
$ws.Range('D41').Value = '12.53'
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('D42').Value = '0.2214'
$ws.Range('E42').Value = '  -0.92%  '
$ws.Range('D43').Value = '0.6833'
$ws.Range('E43').Value = '  +2.67%  '
$ws.Range('D44').Value = '14.29'
$ws.Range('E44').Value = '  +1.12%  '
$ws.Range('D45').Value = '2.346'
$ws.Range('E45').Value = '  +2.28%  '
$ws.Range('D46').Value = '1.002'
$ws.Range('E46').Value = '  -0.35%  '
$ws.Range('D47').Value = '1.368'
$ws.Range('E47').Value = '  +17.67%  '
$ws.Range('D48').Value = '3.637'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('D49').Value = '0.00000000348'
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('E50').Value = '  +9.15%  '
$ws.Range('D51').Value = '1.216'
$ws.Range('E51').Value = '  -0.21%  '
